# edit.ps1
# Applies the OOXML diff changes to before.xlsx via Excel COM interop

$wb = $excel.ActiveWorkbook
$wsCalc = $wb.Worksheets.Item("Calc")
$wsResults = $wb.Worksheets.Item("Results")

# --- Column width changes on "Calc" sheet (columns AX=50, AY=51) ---
$wsCalc.Columns.Item(50).ColumnWidth = 22.7109375
$wsCalc.Columns.Item(51).ColumnWidth = 19.7109375

# --- Column width change on "Results" sheet (column D=4) ---
$wsResults.Columns.Item(4).ColumnWidth = 20.7109375

# --- "Results" sheet header cell C2: unit label change ---
$wsResults.Range("C2").Value = "(ng/g)"

# --- "Calc" sheet cell value changes ---
# Row 3
$wsCalc.Range("AT3").Value = 0.75
$wsCalc.Range("AU3").Value = 0.375

# Row 4
$wsCalc.Range("AT4").Value = 0.75
$wsCalc.Range("AU4").Value = 0.375
$wsCalc.Range("AV4").Value = 0.1218
$wsCalc.Range("AW4").Value = 0.0026896708
$wsCalc.Range("AX4").Value = 0.002695172695222103
$wsCalc.Range("AY4").Value = 2.208268308702792
$wsCalc.Range("BB4").Value = 145.9565715687178
$wsCalc.Range("BC4").Value = 0.5485212227191195
$wsCalc.Range("BD4").Value = 63.80000000000001
$wsCalc.Range("BE4").Value = 1.347586347611051
$wsCalc.Range("BF4").Value = 2.212785464057555

# Row 5
$wsCalc.Range("AT5").Value = 0.75
$wsCalc.Range("AU5").Value = 0.375

# Row 6
$wsCalc.Range("AT6").Value = 0.75
$wsCalc.Range("AU6").Value = 0.375
$wsCalc.Range("AV6").Value = 0.1159
$wsCalc.Range("AW6").Value = 0.0023389926
$wsCalc.Range("AX6").Value = 0.002353196968169543
$wsCalc.Range("AY6").Value = 2.018112683347713
$wsCalc.Range("BB6").Value = 144.5218771154652
$wsCalc.Range("BC6").Value = 0.4967831862381769
$wsCalc.Range("BD6").Value = 57.90000000000001
$wsCalc.Range("BE6").Value = 1.176598484084771
$wsCalc.Range("BF6").Value = 2.030368393588907

# Row 7
$wsCalc.Range("AT7").Value = 0.75
$wsCalc.Range("AU7").Value = 0.375
$wsCalc.Range("AV7").Value = 643.7896
$wsCalc.Range("AX7").Value = 134.1991978658585
$wsCalc.Range("BB7").Value = 6.217884251111586
$wsCalc.Range("BD7").Value = 643731.6
$wsCalc.Range("BE7").Value = 67099.59893292925
$wsCalc.Range("BF7").Value = 20.84519505531908

# Row 8
$wsCalc.Range("AT8").Value = 0.75
$wsCalc.Range("AU8").Value = 0.375
$wsCalc.Range("AV8").Value = 0.1308
$wsCalc.Range("AW8").Value = 0.0025624011
$wsCalc.Range("AX8").Value = 0.002561699555639462
$wsCalc.Range("AY8").Value = 1.959022247706422
$wsCalc.Range("BB8").Value = 144.5736154875932
$wsCalc.Range("BC8").Value = 0.5466913324908758
$wsCalc.Range("BD8").Value = 72.80000000000001
$wsCalc.Range("BE8").Value = 1.280849777819731
$wsCalc.Range("BF8").Value = 1.958485898806928

# Row 9
$wsCalc.Range("AT9").Value = 0.75
$wsCalc.Range("AU9").Value = 0.375
$wsCalc.Range("AV9").Value = 677.1035
$wsCalc.Range("AX9").Value = 163.4758856607667
$wsCalc.Range("BB9").Value = -5.429102174865604
$wsCalc.Range("BD9").Value = 677045.5
$wsCalc.Range("BE9").Value = 81737.94283038336
$wsCalc.Range("BF9").Value = 24.14341170305082

# Row 10
$wsCalc.Range("AT10").Value = 0.75
$wsCalc.Range("AU10").Value = 0.375
$wsCalc.Range("AV10").Value = 0.1141
$wsCalc.Range("AW10").Value = 0.0025474073
$wsCalc.Range("AX10").Value = 0.002517656413905989
$wsCalc.Range("AY10").Value = 2.232609377738826
$wsCalc.Range("BB10").Value = 145.236673811475
$wsCalc.Range("BC10").Value = 0.3792365786174705
$wsCalc.Range("BD10").Value = 56.09999999999999
$wsCalc.Range("BE10").Value = 1.258828206952995
$wsCalc.Range("BF10").Value = 2.206534981512699

# Row 11
$wsCalc.Range("AT11").Value = 0.75
$wsCalc.Range("AU11").Value = 0.375
$wsCalc.Range("AV11").Value = 619.2592
$wsCalc.Range("AX11").Value = 117.4419798265753
$wsCalc.Range("BB11").Value = -4.573091574634732
$wsCalc.Range("BD11").Value = 619201.2
$wsCalc.Range("BE11").Value = 58720.98991328765
$wsCalc.Range("BF11").Value = 18.96491482509672

# Row 12
$wsCalc.Range("AT12").Value = 0.75
$wsCalc.Range("AU12").Value = 0.375
$wsCalc.Range("AV12").Value = 0.1464
$wsCalc.Range("AW12").Value = 0.0033849942
$wsCalc.Range("AX12").Value = 0.003341125105867086
$wsCalc.Range("AY12").Value = 2.312154508196722
$wsCalc.Range("BB12").Value = 145.9747362666762
$wsCalc.Range("BC12").Value = 0.487891868488217
$wsCalc.Range("BD12").Value = 88.4
$wsCalc.Range("BE12").Value = 1.670562552933543
$wsCalc.Range("BF12").Value = 2.282189279963856

# Row 13
$wsCalc.Range("AT13").Value = 0.75
$wsCalc.Range("AU13").Value = 0.375
$wsCalc.Range("AV13").Value = 585.0441
$wsCalc.Range("AX13").Value = 78.41141306457126
$wsCalc.Range("BB13").Value = -1.985162388255983
$wsCalc.Range("BD13").Value = 584986.1
$wsCalc.Range("BE13").Value = 39205.70653228563
$wsCalc.Range("BF13").Value = 13.40264999930283

# --- "Results" sheet cell value changes ---
# Row 3
$wsResults.Range("C3").Value = 13700.97500561748
$wsResults.Range("D3").Value = 1.458374275385301

# Row 4
$wsResults.Range("C4").Value = 2882.82512949779
$wsResults.Range("D4").Value = 0.161942866384525
$wsResults.Range("O4").Value = 0.1218
$wsResults.Range("P4").Value = 0.0026896708
$wsResults.Range("Q4").Value = 145.9565715687178
$wsResults.Range("R4").Value = 0.5485212227191195

# Row 5
$wsResults.Range("C5").Value = 13701.33361959352
$wsResults.Range("D5").Value = 1.718790567296949

# Row 6
$wsResults.Range("C6").Value = 2894.21828608744
$wsResults.Range("D6").Value = 0.17233764492007
$wsResults.Range("O6").Value = 0.1159
$wsResults.Range("P6").Value = 0.0023389926
$wsResults.Range("Q6").Value = 144.5218771154652
$wsResults.Range("R6").Value = 0.4967831862381769

# Row 7
$wsResults.Range("C7").Value = 13703.41462188276
$wsResults.Range("D7").Value = 1.771568459935945
$wsResults.Range("O7").Value = 643.7896
$wsResults.Range("Q7").Value = 6.217884251111586

# Row 8
$wsResults.Range("C8").Value = 2755.624061240333
$wsResults.Range("D8").Value = 0.1972298652588728
$wsResults.Range("O8").Value = 0.1308
$wsResults.Range("P8").Value = 0.0025624011
$wsResults.Range("Q8").Value = 144.5736154875932
$wsResults.Range("R8").Value = 0.5466913324908758

# Row 9
$wsResults.Range("C9").Value = 13702.25911091505
$wsResults.Range("D9").Value = 1.644682893535193
$wsResults.Range("O9").Value = 677.1035
$wsResults.Range("Q9").Value = -5.429102174865604

# Row 10
$wsResults.Range("C10").Value = 2853.925097618624
$wsResults.Range("D10").Value = 0.1451901838312447
$wsResults.Range("O10").Value = 0.1141
$wsResults.Range("P10").Value = 0.0025474073
$wsResults.Range("Q10").Value = 145.236673811475
$wsResults.Range("R10").Value = 0.3792365786174705

# Row 11
$wsResults.Range("C11").Value = 13697.24561754898
$wsResults.Range("D11").Value = 1.369551974574979
$wsResults.Range("O11").Value = 619.2592
$wsResults.Range("Q11").Value = -4.573091574634732

# Row 12
$wsResults.Range("C12").Value = 2734.903355941235
$wsResults.Range("D12").Value = 0.1758825416109763
$wsResults.Range("O12").Value = 0.1464
$wsResults.Range("P12").Value = 0.0033849942
$wsResults.Range("Q12").Value = 145.9747362666762
$wsResults.Range("R12").Value = 0.487891868488217

# Row 13
$wsResults.Range("C13").Value = 13698.12341384287
$wsResults.Range("D13").Value = 1.565772142177249
$wsResults.Range("O13").Value = 585.0441
$wsResults.Range("Q13").Value = -1.985162388255983
